$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Israel Premier League")

    # Row 72 <= data from old row 74
    $ws.Cells.Item(72, 2).Value = 7542640
    $ws.Cells.Item(72, 6).Value = 'MS Ashdod'
    $ws.Cells.Item(72, 7).Value = 'Hapoel Bnei Sakhnin'
    $ws.Cells.Item(72, 8).Value = 0
    $ws.Cells.Item(72, 9).Value = 1
    $ws.Cells.Item(72, 10).Value = 'A'
    $ws.Cells.Item(72, 11).Value = 2.05
    $ws.Cells.Item(72, 12).Value = 3.2
    $ws.Cells.Item(72, 13).Value = 3.5
    $ws.Cells.Item(72, 14).Value = 2.15
    $ws.Cells.Item(72, 15).Value = 3.1
    $ws.Cells.Item(72, 16).Value = 3.2
    $ws.Cells.Item(72, 17).Value = -0.25
    $ws.Cells.Item(72, 18).Value = 1.925
    $ws.Cells.Item(72, 19).Value = 1.925
    $ws.Cells.Item(72, 20).Value = 2.25
    $ws.Cells.Item(72, 21).Value = 1.9
    $ws.Cells.Item(72, 22).Value = 1.95
    $ws.Cells.Item(72, 23).Value = -1
    $ws.Cells.Item(72, 24).Value = -1
    $ws.Cells.Item(72, 25).Value = 2.2
    $ws.Cells.Item(72, 26).Value = -1
    $ws.Cells.Item(72, 27).Value = 0.925
    $ws.Cells.Item(72, 28).Value = -1
    $ws.Cells.Item(72, 29).Value = 0.95

    # Row 73 <= data from old row 72
    $ws.Cells.Item(73, 2).Value = 7542719
    $ws.Cells.Item(73, 6).Value = 'Hapoel Haifa'
    $ws.Cells.Item(73, 7).Value = 'Maccabi Netanya'
    $ws.Cells.Item(73, 8).Value = 2
    $ws.Cells.Item(73, 9).Value = 1
    $ws.Cells.Item(73, 10).Value = 'H'
    $ws.Cells.Item(73, 11).Value = 2.6
    $ws.Cells.Item(73, 12).Value = 3.1
    $ws.Cells.Item(73, 13).Value = 2.6
    $ws.Cells.Item(73, 14).Value = 2.9
    $ws.Cells.Item(73, 15).Value = 3.2
    $ws.Cells.Item(73, 16).Value = 2.3
    $ws.Cells.Item(73, 17).Value = 0.25
    $ws.Cells.Item(73, 18).Value = 1.8
    $ws.Cells.Item(73, 19).Value = 2.05
    $ws.Cells.Item(73, 20).Value = 2.5
    $ws.Cells.Item(73, 21).Value = 2
    $ws.Cells.Item(73, 22).Value = 1.85
    $ws.Cells.Item(73, 23).Value = 1.9
    $ws.Cells.Item(73, 24).Value = -1
    $ws.Cells.Item(73, 25).Value = -1
    $ws.Cells.Item(73, 26).Value = 0.8
    $ws.Cells.Item(73, 27).Value = -1
    $ws.Cells.Item(73, 28).Value = 1
    $ws.Cells.Item(73, 29).Value = -1

    # Row 74 <= data from old row 73
    $ws.Cells.Item(74, 2).Value = 7542639
    $ws.Cells.Item(74, 6).Value = 'Maccabi Bnei Raina'
    $ws.Cells.Item(74, 7).Value = 'Hapoel Jerusalem FC'
    $ws.Cells.Item(74, 8).Value = 1
    $ws.Cells.Item(74, 9).Value = 1
    $ws.Cells.Item(74, 10).Value = 'D'
    $ws.Cells.Item(74, 11).Value = 2.5
    $ws.Cells.Item(74, 12).Value = 3
    $ws.Cells.Item(74, 13).Value = 2.75
    $ws.Cells.Item(74, 14).Value = 2.7
    $ws.Cells.Item(74, 15).Value = 2.8
    $ws.Cells.Item(74, 16).Value = 2.75
    $ws.Cells.Item(74, 17).Value = 0
    $ws.Cells.Item(74, 18).Value = 1.925
    $ws.Cells.Item(74, 19).Value = 1.925
    $ws.Cells.Item(74, 20).Value = 2
    $ws.Cells.Item(74, 21).Value = 2.1
    $ws.Cells.Item(74, 22).Value = 1.775
    $ws.Cells.Item(74, 23).Value = -1
    $ws.Cells.Item(74, 24).Value = 1.8
    $ws.Cells.Item(74, 25).Value = -1
    $ws.Cells.Item(74, 26).Value = 0
    $ws.Cells.Item(74, 27).Value = -0
    $ws.Cells.Item(74, 28).Value = 0
    $ws.Cells.Item(74, 29).Value = -0

    # Row 86 <= data from old row 87
    $ws.Cells.Item(86, 2).Value = 7542727
    $ws.Cells.Item(86, 6).Value = 'Maccabi Bnei Raina'
    $ws.Cells.Item(86, 7).Value = 'Hapoel Bnei Sakhnin'
    $ws.Cells.Item(86, 8).Value = 0
    $ws.Cells.Item(86, 9).Value = 1
    $ws.Cells.Item(86, 10).Value = 'A'
    $ws.Cells.Item(86, 11).Value = 2.1
    $ws.Cells.Item(86, 12).Value = 3.1
    $ws.Cells.Item(86, 13).Value = 3.6
    $ws.Cells.Item(86, 14).Value = 2.45
    $ws.Cells.Item(86, 15).Value = 3
    $ws.Cells.Item(86, 16).Value = 3
    $ws.Cells.Item(86, 17).Value = -0.25
    $ws.Cells.Item(86, 18).Value = 2.075
    $ws.Cells.Item(86, 19).Value = 1.725
    $ws.Cells.Item(86, 20).Value = 2.25
    $ws.Cells.Item(86, 21).Value = 2.05
    $ws.Cells.Item(86, 22).Value = 1.8
    $ws.Cells.Item(86, 23).Value = -1
    $ws.Cells.Item(86, 24).Value = -1
    $ws.Cells.Item(86, 25).Value = 2
    $ws.Cells.Item(86, 26).Value = -1
    $ws.Cells.Item(86, 27).Value = 0.7250000000000001
    $ws.Cells.Item(86, 28).Value = -1
    $ws.Cells.Item(86, 29).Value = 0.8

    # Row 87 <= data from old row 86
    $ws.Cells.Item(87, 2).Value = 7542726
    $ws.Cells.Item(87, 6).Value = 'Hapoel Hadera'
    $ws.Cells.Item(87, 7).Value = 'Maccabi Netanya'
    $ws.Cells.Item(87, 8).Value = 1
    $ws.Cells.Item(87, 9).Value = 4
    $ws.Cells.Item(87, 10).Value = 'A'
    $ws.Cells.Item(87, 11).Value = 3.3
    $ws.Cells.Item(87, 12).Value = 3.5
    $ws.Cells.Item(87, 13).Value = 2
    $ws.Cells.Item(87, 14).Value = 4.333
    $ws.Cells.Item(87, 15).Value = 3.6
    $ws.Cells.Item(87, 16).Value = 1.7
    $ws.Cells.Item(87, 17).Value = 0.75
    $ws.Cells.Item(87, 18).Value = 1.9
    $ws.Cells.Item(87, 19).Value = 1.95
    $ws.Cells.Item(87, 20).Value = 2.5
    $ws.Cells.Item(87, 21).Value = 2
    $ws.Cells.Item(87, 22).Value = 1.85
    $ws.Cells.Item(87, 23).Value = -1
    $ws.Cells.Item(87, 24).Value = -1
    $ws.Cells.Item(87, 25).Value = 0.7
    $ws.Cells.Item(87, 26).Value = -1
    $ws.Cells.Item(87, 27).Value = 0.95
    $ws.Cells.Item(87, 28).Value = 1
    $ws.Cells.Item(87, 29).Value = -1

    # Row 108 <= data from old row 110
    $ws.Cells.Item(108, 2).Value = 7542737
    $ws.Cells.Item(108, 6).Value = 'MS Ashdod'
    $ws.Cells.Item(108, 7).Value = 'Hapoel Haifa'
    $ws.Cells.Item(108, 8).Value = 0
    $ws.Cells.Item(108, 9).Value = 1
    $ws.Cells.Item(108, 10).Value = 'A'
    $ws.Cells.Item(108, 11).Value = 3
    $ws.Cells.Item(108, 12).Value = 3.2
    $ws.Cells.Item(108, 13).Value = 2.45
    $ws.Cells.Item(108, 14).Value = 3.2
    $ws.Cells.Item(108, 15).Value = 3.25
    $ws.Cells.Item(108, 16).Value = 2.3
    $ws.Cells.Item(108, 17).Value = 0.25
    $ws.Cells.Item(108, 18).Value = 1.85
    $ws.Cells.Item(108, 19).Value = 2
    $ws.Cells.Item(108, 20).Value = 2.25
    $ws.Cells.Item(108, 21).Value = 1.875
    $ws.Cells.Item(108, 22).Value = 1.975
    $ws.Cells.Item(108, 23).Value = -1
    $ws.Cells.Item(108, 24).Value = -1
    $ws.Cells.Item(108, 25).Value = 1.3
    $ws.Cells.Item(108, 26).Value = -1
    $ws.Cells.Item(108, 27).Value = 1
    $ws.Cells.Item(108, 28).Value = -1
    $ws.Cells.Item(108, 29).Value = 0.9750000000000001

    # Row 109 <= data from old row 108
    $ws.Cells.Item(109, 2).Value = 7542735
    $ws.Cells.Item(109, 6).Value = 'Hapoel Petah Tikva'
    $ws.Cells.Item(109, 7).Value = 'Maccabi Netanya'
    $ws.Cells.Item(109, 8).Value = 2
    $ws.Cells.Item(109, 9).Value = 0
    $ws.Cells.Item(109, 10).Value = 'H'
    $ws.Cells.Item(109, 11).Value = 3.75
    $ws.Cells.Item(109, 12).Value = 3.6
    $ws.Cells.Item(109, 13).Value = 1.909
    $ws.Cells.Item(109, 14).Value = 3.8
    $ws.Cells.Item(109, 15).Value = 3.75
    $ws.Cells.Item(109, 16).Value = 1.85
    $ws.Cells.Item(109, 17).Value = 0.5
    $ws.Cells.Item(109, 18).Value = 1.95
    $ws.Cells.Item(109, 19).Value = 1.9
    $ws.Cells.Item(109, 20).Value = 2.5
    $ws.Cells.Item(109, 21).Value = 1.975
    $ws.Cells.Item(109, 22).Value = 1.875
    $ws.Cells.Item(109, 23).Value = 2.8
    $ws.Cells.Item(109, 24).Value = -1
    $ws.Cells.Item(109, 25).Value = -1
    $ws.Cells.Item(109, 26).Value = 0.95
    $ws.Cells.Item(109, 27).Value = -1
    $ws.Cells.Item(109, 28).Value = -1
    $ws.Cells.Item(109, 29).Value = 0.875

    # Row 110 <= data from old row 109
    $ws.Cells.Item(110, 2).Value = 7542736
    $ws.Cells.Item(110, 6).Value = 'Hapoel Jerusalem FC'
    $ws.Cells.Item(110, 7).Value = 'Hapoel Bnei Sakhnin'
    $ws.Cells.Item(110, 8).Value = 0
    $ws.Cells.Item(110, 9).Value = 0
    $ws.Cells.Item(110, 10).Value = 'D'
    $ws.Cells.Item(110, 11).Value = 2.2
    $ws.Cells.Item(110, 12).Value = 3.4
    $ws.Cells.Item(110, 13).Value = 3.2
    $ws.Cells.Item(110, 14).Value = 2.375
    $ws.Cells.Item(110, 15).Value = 3.1
    $ws.Cells.Item(110, 16).Value = 3.1
    $ws.Cells.Item(110, 17).Value = -0.25
    $ws.Cells.Item(110, 18).Value = 2.05
    $ws.Cells.Item(110, 19).Value = 1.8
    $ws.Cells.Item(110, 20).Value = 2
    $ws.Cells.Item(110, 21).Value = 1.925
    $ws.Cells.Item(110, 22).Value = 1.925
    $ws.Cells.Item(110, 23).Value = -1
    $ws.Cells.Item(110, 24).Value = 2.1
    $ws.Cells.Item(110, 25).Value = -1
    $ws.Cells.Item(110, 26).Value = -0.5
    $ws.Cells.Item(110, 27).Value = 0.4
    $ws.Cells.Item(110, 28).Value = -1
    $ws.Cells.Item(110, 29).Value = 0.925

    # Row 201: update U and V
    $ws.Cells.Item(201, 21).Value = 2
    $ws.Cells.Item(201, 22).Value = 1.85

    # Row 203: update U and V
    $ws.Cells.Item(203, 21).Value = 1.825
    $ws.Cells.Item(203, 22).Value = 2.025

